$d = $word.ActiveDocument

$replacements = @(
    @("92÷6=", "11÷9="),
    @("41÷3=", "96÷6="),
    @("53÷4=", "23÷6="),
    @("84÷7=", "19÷4="),
    @("38÷3=", "41÷4="),
    @("43÷7=", "88÷2="),
    @("94÷6=", "39÷9="),
    @("71÷4=", "95÷8="),
    @("84÷2=", "38÷3="),
    @("65÷8=", "42÷6="),
    @("73÷9=", "72÷5="),
    @("68÷4=", "89÷5="),
    @("67÷2=", "11÷4="),
    @("35÷7=", "71÷5="),
    @("80÷9=", "68÷2="),
    @("67÷8=", "85÷9="),
    @("88÷7=", "39÷3="),
    @("75÷9=", "71÷8="),
    @("61÷7=", "89÷6="),
    @("96÷5=", "39÷5="),
    @("33÷4=", "34÷9="),
    @("19÷9=", "50÷8="),
    @("74÷7=", "65÷2="),
    @("72÷6=", "54÷3="),
    @("53÷2=", "67÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1)
}
